# Updated CHE_grids model - 2025-08-19 00:07
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) ev_charging_uc sheet: timeslice lookup table (rows 12-14) changes
#    from the old day/night timeslice-list definitions to the new
#    day_night / timeslice / AllSaAllH based definitions.
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ev_charging_uc")
$ws1.Range("B12").Value = "day_night"
$ws1.Range("C12").Value = "timeslice"
$ws1.Range("B13").Value = "D"
$ws1.Range("C13").Value = "AllSaAllH"
$ws1.Range("B14").Value = "N"
$ws1.Range("C14").Value = "AllSaAllH"

# -----------------------------------------------------------------
# 2) Add two new worksheets at the end: re_profiles and load_shapes
# -----------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRe = $wb.Worksheets.Add($null, $lastSheet)
$wsRe.Name = "re_profiles"

$wsLoad = $wb.Worksheets.Add($null, $wsRe)
$wsLoad.Name = "load_shapes"

# --- re_profiles content ---
$wsRe.Range("B2").Value = "~TFM_DINS-AT"
$wsRe.Range("B3").Value = "timeslice"
$wsRe.Range("C3").Value = "com_fr"
$wsRe.Range("D3").Value = "process"
$wsRe.Range("E3").Value = "commodity"
$wsRe.Range("B4").Value = "AllSaAllH"
$wsRe.Range("C4").Value = 0.99999979999999999
$wsRe.Range("D4").Value = "IMPNRGZ"
$wsRe.Range("E4").Value = "elc_spv-CHE"

$wsRe.Range("H2").Value = "~TFM_DINS-AT"
$wsRe.Range("H3").Value = "timeslice"
$wsRe.Range("I3").Value = "com_fr"
$wsRe.Range("J3").Value = "process"
$wsRe.Range("K3").Value = "commodity"
$wsRe.Range("H4").Value = "AllSaAllH"
$wsRe.Range("I4").Value = 0.999999999999855
$wsRe.Range("J4").Value = "IMPNRGZ"
$wsRe.Range("K4").Value = "elc_won-CHE"

$wsRe.Range("M2").Value = "~TFM_INS-AT"
$wsRe.Range("M3").Value = "timeslice"
$wsRe.Range("N3").Value = "ncap_afs"
$wsRe.Range("O3").Value = "pset_ci"
$wsRe.Range("M4").Value = "AllS"
$wsRe.Range("N4").Value = 1.2
$wsRe.Range("O4").Value = "hydro"

# --- load_shapes content ---
$wsLoad.Range("B2").Value = "~TFM_DINS-AT"
$wsLoad.Range("B3").Value = "g_yrfr"
$wsLoad.Range("C3").Value = "com_fr"
$wsLoad.Range("D3").Value = "timeslice"
$wsLoad.Range("E3").Value = "commodity"
$wsLoad.Range("B4").Value = 1
$wsLoad.Range("C4").Value = 1.0000000000000002
$wsLoad.Range("D4").Value = "AllSaAllH"
$wsLoad.Range("E4").Value = "elc_roadtransport"

$wsLoad.Range("H2").Value = "~TFM_DINS-AT"
$wsLoad.Range("H3").Value = "commodity"
$wsLoad.Range("I3").Value = "timeslice"
$wsLoad.Range("J3").Value = "com_fr"
$wsLoad.Range("H4").Value = "elc_buildings"
$wsLoad.Range("I4").Value = "AllSaAllH"
$wsLoad.Range("J4").Value = 1
$wsLoad.Range("H5").Value = "elc_industry"
$wsLoad.Range("I5").Value = "AllSaAllH"
$wsLoad.Range("J5").Value = 1

$wsLoad.Range("M2").Value = "~TFM_DINS-AT"
$wsLoad.Range("M3").Value = "commodity"
$wsLoad.Range("N3").Value = "timeslice"
$wsLoad.Range("O3").Value = "com_pkflx"
$wsLoad.Range("M4").Value = "ELC"
$wsLoad.Range("N4").Value = "AllSaAllH"
$wsLoad.Range("O4").Value = 0.29211718079128235

# -----------------------------------------------------------------
# 3) Window view: land on the new last sheet (load_shapes) as active
# -----------------------------------------------------------------
$wsLoad.Activate()
